$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-01-26 08:57:36"
$wsZhCn.Range("G4").Value = "2016-01-26 08:58:23"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-01-26 08:57:47"
$wsDeDe.Range("G4").Value = "2016-01-26 08:58:41"
